$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the daily/weekly logic sheet.
# It is inserted as a new row 118, pushing the existing rows 118-152 down
# to rows 119-153 (all their contents, including styles, stay identical -
# only the row number shifts).
$ws.Rows.Item(118).Insert()

# Populate the newly inserted row 118 with the new record's data.
$ws.Cells.Item(118, 1).Value = 11
$ws.Cells.Item(118, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(118, 3).Value = "Bíobío"
$ws.Cells.Item(118, 4).Value = 44988
$ws.Cells.Item(118, 5).Value = 8
$ws.Cells.Item(118, 6).Value = "Fruta"
$ws.Cells.Item(118, 7).Value = 100108
$ws.Cells.Item(118, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(118, 9).Value = 100108002
$ws.Cells.Item(118, 10).Value = "Mango"
$ws.Cells.Item(118, 11).Value = "Sin especificar"
$ws.Cells.Item(118, 12).Value = "Primera"
$ws.Cells.Item(118, 13).Value = 200
$ws.Cells.Item(118, 14).Value = 7000
$ws.Cells.Item(118, 15).Value = 8000
$ws.Cells.Item(118, 16).Value = 7500
$ws.Cells.Item(118, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(118, 18).Value = "Perú"
$ws.Cells.Item(118, 19).Value = 1875
$ws.Cells.Item(118, 20).Value = 4
